# Insert a new price-record row at row 210 of the single data sheet,
# pushing the former rows 210-291 down to 211-292 (dimension grows from
# A1:R291 to A1:R292). The new row carries a new weekly observation for
# "Feria Lagunitas de Puerto Montt" / Acelga.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 210 (and everything below it) down by one row while
# keeping formatting (date style on column D, etc.) intact.
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(210, 1).Value = 4
$ws.Cells.Item(210, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(210, 3).Value = "Los Lagos"
$ws.Cells.Item(210, 4).Value = 45007
$ws.Cells.Item(210, 5).Value = 10
$ws.Cells.Item(210, 6).Value = 100112009
$ws.Cells.Item(210, 7).Value = "Acelga"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 15
$ws.Cells.Item(210, 11).Value = 10000
$ws.Cells.Item(210, 12).Value = 10000
$ws.Cells.Item(210, 13).Value = 10000
$ws.Cells.Item(210, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(210, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(210, 16).Value = 833
$ws.Cells.Item(210, 17).Value = 12
$ws.Cells.Item(210, 18).Value = "Hortaliza"
